$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.837.82'
$ws.Range("E2").Value = '  -4.62%  '

$ws.Range("D3").Value = '2.691.66'
$ws.Range("E3").Value = '  -7.52%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").Value = '''495.93'
$ws.Range("E5").Value = '  -6.06%  '

$ws.Range("D6").Value = '''134.99'
$ws.Range("E6").Value = '  -5.51%  '

$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  -0.19%  '

$ws.Range("D8").Value = '''0.525'
$ws.Range("E8").Value = '  -5.04%  '

$ws.Range("D9").Value = '2.698.41'
$ws.Range("E9").Value = '  -7.30%  '

$ws.Range("B10").Value = 'Toncoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D10").Value = '''5.86'
$ws.Range("E10").Value = '  +0.19%  '

$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").Value = '''0.101'
$ws.Range("E11").Value = '  -5.71%  '

$ws.Range("D12").Value = '''0.340'
$ws.Range("E12").Value = '  -3.21%  '

$ws.Range("E13").Value = '  +1.08%  '

$ws.Range("D14").Value = '3.166.59'
$ws.Range("E14").Value = '  -7.23%  '

$ws.Range("D15").Value = '57.947.11'
$ws.Range("E15").Value = '  -4.59%  '

$ws.Range("D16").Value = '''21.02'
$ws.Range("E16").Value = '  -6.93%  '

$ws.Range("D17").Value = '2.691.21'
$ws.Range("E17").Value = '  -7.36%  '

$ws.Range("D18").Value = '''0.0000133'
$ws.Range("E18").Value = '  -5.35%  '

$ws.Range("D19").Value = '''4.63'
$ws.Range("E19").Value = '  -5.80%  '

$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").Value = '''10.72'
$ws.Range("E20").Value = '  -6.93%  '

$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").Value = '''339.86'
$ws.Range("E21").Value = '  -5.60%  '

$ws.Range("D22").Value = '''6.11'
$ws.Range("E22").Value = '  -6.36%  '

$ws.Range("D23").Value = '''0.998'
$ws.Range("E23").Value = '  -0.15%  '

$ws.Range("E24").Value = '  -1.20%  '

$ws.Range("D25").Value = '''61.47'
$ws.Range("E25").Value = '  -3.04%  '

$ws.Range("D26").Value = '''0.418'
$ws.Range("E26").Value = '  -7.27%  '

$ws.Range("D27").Value = '''0.169'
$ws.Range("E27").Value = '  -5.61%  '

$ws.Range("D28").Value = '''0.999'
$ws.Range("E28").Value = '  -0.10%  '

$ws.Range("D29").Value = '''7.29'
$ws.Range("E29").Value = '  -4.73%  '

$ws.Range("D30").Value = '0.0₃0805'
$ws.Range("E30").Value = '  -6.00%  '

$ws.Range("D31").Value = '''0.999'
$ws.Range("E31").Value = '  -0.11%  '

$ws.Range("E32").Value = '  -5.58%  '

$ws.Range("D33").Value = '''18.78'
$ws.Range("E33").Value = '  -4.35%  '

$ws.Range("D34").Value = '''146.90'
$ws.Range("E34").Value = '  -3.29%  '

$ws.Range("D35").Value = '''4.11'
$ws.Range("E35").Value = '  -5.18%  '

$ws.Range("D36").Value = '''5.21'
$ws.Range("E36").Value = '  -6.28%  '

$ws.Range("D37").Value = '''0.899'
$ws.Range("E37").Value = '  -10.09%  '

$ws.Range("D38").Value = '''1.11'
$ws.Range("E38").Value = '  -7.51%  '

$ws.Range("D39").Value = '''36.06'
$ws.Range("E39").Value = '  -4.91%  '

$ws.Range("B40").Value = 'Filecoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D40").Value = '''3.47'
$ws.Range("E40").Value = '  -5.63%  '

$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D41").Value = '''0.997'
$ws.Range("E41").Value = '  -0.19%  '

$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '2.135.92'
$ws.Range("E42").Value = '  -8.32%  '

$ws.Range("D43").Value = '''1.34'
$ws.Range("E43").Value = '  -8.69%  '

$ws.Range("D44").Value = '''0.0546'
$ws.Range("E44").Value = '  -4.10%  '

$ws.Range("D45").Value = '''0.589'
$ws.Range("E45").Value = '  -8.35%  '

$ws.Range("D46").Value = '''10.32'
$ws.Range("E46").Value = '  -0.03%  '

$ws.Range("D47").Value = '''18.69'
$ws.Range("E47").Value = '  -10.03%  '

$ws.Range("D48").Value = '''0.0222'
$ws.Range("E48").Value = '  -4.33%  '

$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").Value = '''0.0878'
$ws.Range("E49").Value = '  -4.94%  '

$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").Value = '''4.49'
$ws.Range("E50").Value = '  -7.44%  '

$ws.Range("D51").Value = '''17.29'
$ws.Range("E51").Value = '  -5.64%  '
